$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 56, shifting rows 56:112 down to 57:113
$ws.Rows.Item(56).Insert(-4121)

# Populate the new row 56 with the new data record.
# Column layout: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg, F Tipo,
# G Producto ID, H Producto, I Categoria ID, J Categoria, K Variedad,
# L Calidad, M Volumen, N Precio minimo, O Precio maximo,
# P Precio promedio ponderado, Q Unidad de comercializacion, R Origen,
# S Precio $/Kg, T Kg / unidad
$ws.Cells.Item(56, 1).Value = 1
$ws.Cells.Item(56, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(56, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(56, 4).Value = 44902
$ws.Cells.Item(56, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(56, 5).Value = 15
$ws.Cells.Item(56, 6).Value = "Fruta"
$ws.Cells.Item(56, 7).Value = 100102
$ws.Cells.Item(56, 8).Value = "Cítricos"
$ws.Cells.Item(56, 9).Value = 100102005
$ws.Cells.Item(56, 10).Value = "Naranja"
$ws.Cells.Item(56, 11).Value = "Lane Late"
$ws.Cells.Item(56, 12).Value = "Segunda"
$ws.Cells.Item(56, 13).Value = 300
$ws.Cells.Item(56, 14).Value = 950
$ws.Cells.Item(56, 15).Value = 1000
$ws.Cells.Item(56, 16).Value = 975
$ws.Cells.Item(56, 17).Value = "$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(56, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(56, 19).Value = 975
$ws.Cells.Item(56, 20).Value = 1
